$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 139
$ws.Range("I2").Value = 309
$ws.Range("J2").Value = 1495
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 425
$ws.Range("M2").Value = 22
$ws.Range("N2").Value = 257
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 19
$ws.Range("S2").Value = 173
$ws.Range("T2").Value = 234
$ws.Range("U2").Value = 24
$ws.Range("V2").Value = 2226
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 2241
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 36
$ws.Range("AA2").Value = 28
